$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: consolidate "A" + " " + "slide" into a single run "A slide".
# The text already reads as "A slide", so re-assign via a distinct
# intermediate value first to force the writer to rebuild the run list.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "A slide__tmp__"
$titleRange.Text = "A slide"

# Table cell (row 1, col 2): consolidate "a" + " " + "table" into "a table".
$table = $s.Shapes.Item(3).Table
$cellRange = $table.Rows.Item(1).Cells.Item(2).Shape.TextFrame.TextRange
$cellRange.Text = "a table__tmp__"
$cellRange.Text = "a table"
